# Auto-save via app Streamlit
# One of the duplicated "Claudine Fleury" booking rows (row 51) is removed.
# This shifts the TOTAL row (previously row 52) up to row 51, and the
# worksheet's used range shrinks from A1:O52 to A1:O51.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(51).Delete()
